$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 272, shifting rows 272:291 down to 273:292
$ws.Rows.Item(272).Insert()

# Populate the new row 272 with the new weekly data point.
$ws.Cells.Item(272, 1).Value = 5
$ws.Cells.Item(272, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(272, 3).Value = "Maule"
$ws.Cells.Item(272, 4).Value = 44783
$ws.Cells.Item(272, 5).Value = 7
$ws.Cells.Item(272, 6).Value = 100112009
$ws.Cells.Item(272, 7).Value = "Acelga"
$ws.Cells.Item(272, 8).Value = "Sin especificar"
$ws.Cells.Item(272, 9).Value = "Primera"
$ws.Cells.Item(272, 10).Value = 500
$ws.Cells.Item(272, 11).Value = 3000
$ws.Cells.Item(272, 12).Value = 3000
$ws.Cells.Item(272, 13).Value = 3000
$ws.Cells.Item(272, 14).Value = "`$/docena de atados (4 kilos)"
$ws.Cells.Item(272, 15).Value = "Región del Maule"
$ws.Cells.Item(272, 16).Value = 750
$ws.Cells.Item(272, 17).Value = 4
$ws.Cells.Item(272, 18).Value = "Hortaliza"

# Copy the date cell style from the row above (style index 2, date number format)
$ws.Cells.Item(271, 4).Copy()
$ws.Cells.Item(272, 4).PasteSpecial(-4122)
$excel.CutCopyMode = $false
